$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 11
$ws.Range("H11").Value2 = 97.29412000000001
$ws.Range("I11").Value2 = 97.29412000000001
$ws.Range("K11").Value2 = 97.29412000000001
$ws.Range("M11").Value2 = 42.70587999999999
# row 55
$ws.Range("H55").Value2 = 251.16667
$ws.Range("J55").Value2 = 313.7143
$ws.Range("L55").Value2 = 313.7143
$ws.Range("N55").Value2 = -741.7143
# row 92
$ws.Range("H92").Value2 = 1233.75
$ws.Range("J92").Value2 = 1093.8
$ws.Range("L92").Value2 = 1093.8
$ws.Range("N92").Value2 = -3589.8
# row 111
$ws.Range("H111").Value2 = 1838.8055
$ws.Range("I111").Value2 = 5313.857
$ws.Range("K111").Value2 = 15941.571
$ws.Range("M111").Value2 = -12874.571
# row 132
$ws.Range("H132").Value2 = 8348.833000000001
$ws.Range("I132").Value2 = 10452.786
$ws.Range("K132").Value2 = 31358.358
$ws.Range("M132").Value2 = -28828.358
# row 138
$ws.Range("H138").Value2 = 2540.9375
$ws.Range("J138").Value2 = 2553.5933
$ws.Range("L138").Value2 = 7660.7799
$ws.Range("N138").Value2 = -17940.7799

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value2 = 2428.4
$ws.Range("I32").Value2 = 2411.647
$ws.Range("K32").Value2 = 2411.647
$ws.Range("M32").Value2 = -2124.647
# row 61
$ws.Range("H61").Value2 = 3101
$ws.Range("I61").Value2 = 2961.6
$ws.Range("J61").Value2 = 3333.3333
$ws.Range("K61").Value2 = 2961.6
$ws.Range("L61").Value2 = 3333.3333
$ws.Range("M61").Value2 = -2749.6
$ws.Range("N61").Value2 = -3757.3333
# row 86
$ws.Range("H86").Value2 = 90000
$ws.Range("J86").Value2 = 90000
$ws.Range("L86").Value2 = 90000
$ws.Range("N86").Value2 = -92372
# row 89
$ws.Range("H89").Value2 = 90000
$ws.Range("J89").Value2 = 90000
$ws.Range("L89").Value2 = 270000
$ws.Range("N89").Value2 = -281856
# row 113
$ws.Range("H113").Value2 = 0
$ws.Range("J113").Value2 = 0
$ws.Range("L113").Value2 = 0
$ws.Range("N113").ClearContents()
# row 132
$ws.Range("H132").Value2 = 0
$ws.Range("I132").Value2 = 0
$ws.Range("K132").Value2 = 0
$ws.Range("M132").ClearContents()
# row 136
$ws.Range("H136").Value2 = 3101
$ws.Range("I136").Value2 = 2961.6
$ws.Range("J136").Value2 = 3333.3333
$ws.Range("K136").Value2 = 8884.799999999999
$ws.Range("L136").Value2 = 9999.999899999999
$ws.Range("M136").Value2 = -6334.799999999999
$ws.Range("N136").Value2 = -15099.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 97
$ws.Range("H97").Value2 = 30000
$ws.Range("J97").Value2 = 30000
$ws.Range("L97").Value2 = 30000
$ws.Range("N97").Value2 = -31982
# row 99
$ws.Range("H99").Value2 = 57825.555
$ws.Range("I99").Value2 = 68677.39999999999
$ws.Range("K99").Value2 = 68677.39999999999
$ws.Range("M99").Value2 = -67179.39999999999
# row 134
$ws.Range("H134").Value2 = 5003.5
$ws.Range("I134").Value2 = 2000
$ws.Range("J134").Value2 = 6004.6665
$ws.Range("K134").Value2 = 6000
$ws.Range("L134").Value2 = 18013.9995
$ws.Range("M134").Value2 = -3465
$ws.Range("N134").Value2 = -23083.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 23
$ws.Range("H23").Value2 = 3343336.2
$ws.Range("I23").Value2 = 3343336.2
$ws.Range("K23").Value2 = 3343336.2
$ws.Range("M23").Value2 = -3343096.2
# row 27
$ws.Range("H27").Value2 = 3343336.2
$ws.Range("I27").Value2 = 3343336.2
$ws.Range("K27").Value2 = 3343336.2
$ws.Range("M27").Value2 = -3343144.2
# row 43
$ws.Range("H43").Value2 = 37999.668
$ws.Range("J43").Value2 = 37999.668
$ws.Range("L43").Value2 = 37999.668
$ws.Range("N43").Value2 = -38367.668
# row 101
$ws.Range("H101").Value2 = 37999.668
$ws.Range("J101").Value2 = 37999.668
$ws.Range("L101").Value2 = 37999.668
$ws.Range("N101").Value2 = -44489.668
# row 119
$ws.Range("H119").Value2 = 130490
$ws.Range("J119").Value2 = 130490
$ws.Range("L119").Value2 = 130490
$ws.Range("N119").Value2 = -140166
# row 134
$ws.Range("H134").Value2 = 4466.3335
$ws.Range("I134").Value2 = 4774.75
$ws.Range("K134").Value2 = 14324.25
$ws.Range("M134").Value2 = -11789.25
# row 141
$ws.Range("H141").Value2 = 123145
$ws.Range("I141").Value2 = 62997
$ws.Range("J141").Value2 = 133169.67
$ws.Range("K141").Value2 = 62997
$ws.Range("L141").Value2 = 133169.67
$ws.Range("M141").Value2 = -57817
$ws.Range("N141").Value2 = -143529.67

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 109
$ws.Range("H109").Value2 = 3802
$ws.Range("I109").Value2 = 2805.6667
$ws.Range("J109").Value2 = 4345.4546
$ws.Range("K109").Value2 = 8417.000100000001
$ws.Range("L109").Value2 = 13036.3638
$ws.Range("M109").Value2 = -7377.000100000001
$ws.Range("N109").Value2 = -15116.3638
# row 112
$ws.Range("H112").Value2 = 2745.5
$ws.Range("I112").Value2 = 2010
$ws.Range("J112").Value2 = 3481
$ws.Range("K112").Value2 = 6030
$ws.Range("L112").Value2 = 10443
$ws.Range("M112").Value2 = -4922
$ws.Range("N112").Value2 = -12659
# row 113
$ws.Range("H113").Value2 = 1800.375
$ws.Range("J113").Value2 = 1800.375
$ws.Range("L113").Value2 = 5401.125
$ws.Range("N113").Value2 = -9741.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value2 = 7602.25
$ws.Range("J70").Value2 = 7602.25
$ws.Range("L70").Value2 = 7602.25
$ws.Range("N70").Value2 = -8142.25
# row 73
$ws.Range("H73").Value2 = 7602.25
$ws.Range("J73").Value2 = 7602.25
$ws.Range("L73").Value2 = 7602.25
$ws.Range("N73").Value2 = -9474.25
# row 97
$ws.Range("H97").Value2 = 1608.25
$ws.Range("I97").Value2 = 1695.1428
$ws.Range("K97").Value2 = 1695.1428
$ws.Range("M97").Value2 = -1199.1428
# row 100
$ws.Range("H100").Value2 = 49994
$ws.Range("J100").Value2 = 49994
$ws.Range("L100").Value2 = 49994
$ws.Range("N100").Value2 = -52158
# row 102
$ws.Range("H102").Value2 = 9216.286
$ws.Range("I102").Value2 = 0
$ws.Range("J102").Value2 = 9216.286
$ws.Range("K102").Value2 = 0
$ws.Range("L102").Value2 = 9216.286
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value2 = -12460.286
# row 122
$ws.Range("H122").Value2 = 9619096
$ws.Range("I122").Value2 = 10991824
$ws.Range("J122").Value2 = 10000
$ws.Range("K122").Value2 = 32975472
$ws.Range("L122").Value2 = 30000
$ws.Range("M122").Value2 = -32973022
$ws.Range("N122").Value2 = -34900
# row 132
$ws.Range("H132").Value2 = 4749.5
$ws.Range("I132").Value2 = 5999
$ws.Range("K132").Value2 = 17997
$ws.Range("M132").Value2 = -15467
# row 135
$ws.Range("H135").Value2 = 69979.09
$ws.Range("J135").Value2 = 69979.09
$ws.Range("L135").Value2 = 69979.09
$ws.Range("N135").Value2 = -80119.09

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value2 = 2229.875
# row 43
$ws.Range("H43").Value2 = 15000
$ws.Range("I43").Value2 = 0
$ws.Range("K43").Value2 = 0
$ws.Range("M43").ClearContents()
# row 61
$ws.Range("H61").Value2 = 3004.2903
$ws.Range("I61").Value2 = 2909.64
$ws.Range("K61").Value2 = 2909.64
$ws.Range("M61").Value2 = -2707.64
# row 100
$ws.Range("H100").Value2 = 3871.375
$ws.Range("I100").Value2 = 1995
$ws.Range("J100").Value2 = 6998.6665
$ws.Range("K100").Value2 = 1995
$ws.Range("L100").Value2 = 6998.6665
$ws.Range("M100").Value2 = -1454
$ws.Range("N100").Value2 = -8080.6665
# row 103
$ws.Range("H103").Value2 = 24999.666
$ws.Range("J103").Value2 = 24999.666
$ws.Range("L103").Value2 = 24999.666
$ws.Range("N103").Value2 = -27343.666
# row 106
$ws.Range("H106").Value2 = 20312.334
$ws.Range("J106").Value2 = 20312.334
$ws.Range("L106").Value2 = 20312.334
$ws.Range("N106").Value2 = -22836.334
# row 113
$ws.Range("H113").Value2 = 3004.2903
$ws.Range("I113").Value2 = 2909.64
$ws.Range("K113").Value2 = 2909.64
$ws.Range("M113").Value2 = -739.6399999999999
# row 126
$ws.Range("H126").Value2 = 2229.875
# row 136
$ws.Range("H136").Value2 = 2790.641
$ws.Range("I136").Value2 = 2205.1667
$ws.Range("K136").Value2 = 6615.500100000001
$ws.Range("M136").Value2 = -4065.500100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 21
$ws.Range("H21").Value2 = 0
$ws.Range("I21").Value2 = 0
$ws.Range("K21").Value2 = 0
$ws.Range("M21").ClearContents()
# row 29
$ws.Range("H29").Value2 = 16670
$ws.Range("I29").Value2 = 22505
$ws.Range("J29").Value2 = 5000
$ws.Range("K29").Value2 = 22505
$ws.Range("L29").Value2 = 5000
$ws.Range("M29").Value2 = -22215
$ws.Range("N29").Value2 = -5580
# row 35
$ws.Range("H35").Value2 = 0
$ws.Range("I35").Value2 = 0
$ws.Range("K35").Value2 = 0
$ws.Range("M35").ClearContents()
# row 128
$ws.Range("H128").Value2 = 0
$ws.Range("J128").Value2 = 0
$ws.Range("L128").Value2 = 0
$ws.Range("N128").ClearContents()
